$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 9 (the standalone "SHANMUGAPRIYAN S" row) - shifts rows 10-19 up to 9-18
$ws.Rows("9:9").Delete()

# Fix the "Mr. BALU M S" label - remove the stray leading space
$ws.Range("A4").Value = "Mr. BALU M S"

# Update the active selection
$ws.Range("C23").Select() | Out-Null
